# Update odds data on the active worksheet to reflect the latest FlashScore snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (Deportes Iquique vs Everton) ---
$ws.Range("G2").Value  = 2.63
$ws.Range("H2").Value  = 3.2
$ws.Range("J2").Value  = 1.06
$ws.Range("K2").Value  = 9.5
$ws.Range("L2").Value  = 1.3
$ws.Range("M2").Value  = 3.4
$ws.Range("N2").Value  = 2
$ws.Range("O2").Value  = 1.8
$ws.Range("P2").Value  = 1.4
$ws.Range("Q2").Value  = 2.75
$ws.Range("R2").Value  = 1.73
$ws.Range("S2").Value  = 2
$ws.Range("T2").Value  = 9
$ws.Range("V2").Value  = 10
$ws.Range("X2").Value  = 21
$ws.Range("Y2").Value  = 29
$ws.Range("Z2").Value  = 9.5
$ws.Range("AB2").Value = 13
$ws.Range("AC2").Value = 41
$ws.Range("AD2").Value = 201
$ws.Range("AE2").Value = 9
$ws.Range("AI2").Value = 21
$ws.Range("AJ2").Value = 29

# --- Row 3 (Santa Fe vs Once Caldas) ---
$ws.Range("G3").Value = 1.85
$ws.Range("I3").Value = 4.5

# --- Row 6 (Vikingur Reykjavik vs KR Reykjavik) ---
$ws.Range("G6").Value  = 1.55
$ws.Range("H6").Value  = 4.65
$ws.Range("I6").Value  = 4.55
$ws.Range("M6").Value  = 6.7
$ws.Range("T6").Value  = 17.5
$ws.Range("X6").Value  = 11
$ws.Range("Y6").Value  = 14
$ws.Range("AB6").Value = 12.5
$ws.Range("AF6").Value = 40
$ws.Range("AG6").Value = 16.5
$ws.Range("AH6").Value = 80
$ws.Range("AI6").Value = 35
$ws.Range("AJ6").Value = 26

# --- Row 7 (Drogheda vs Shamrock Rovers) ---
$ws.Range("I7").Value = 1.8

# --- Row 8 (St. Patricks vs Shelbourne) ---
$ws.Range("K8").Value = 9
